# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the first file row (1fe9a58a-...) on both the zh-cn and de-de
# report sheets, simulating a fresh handback-status report run.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-24 11:41:55"
$ws_zhcn.Range("H2").Value = "2016-03-24 11:42:22"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-24 11:42:00"
$ws_dede.Range("H2").Value = "2016-03-24 11:42:29"
